$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet so its data/format survive as the
#    (unchanged) new "2022-Q2" sheet, then rename the original to "2022-Q3"
#    and fill it in with the new quarter's numbers.
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Copy it right after itself -> becomes "2022-Q2 (2)" just after the original.
$q2Sheet.Copy($null, $q2Sheet)

# The original (still holding the old Q2 numbers) becomes the brand-new Q3 tab.
$q2Sheet.Name = "2022-Q3"

# The freshly made copy takes over the "2022-Q2" name (data unchanged).
$q2Copy = $wb.Worksheets.Item(3)
$q2Copy.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 2) Write the 2022-Q3 figures into the renamed sheet (rank dropped 6 -> 4,
#    position metrics updated).
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

$q3Sheet.Range("D2").Value = "'0.07"
$q3Sheet.Range("E2").Value = "'86.27"
$q3Sheet.Range("F2").Value = "'3.97"
$q3Sheet.Range("G2").Value = "'0.0028"
$q3Sheet.Range("H2").Value = 4

$q3Sheet.Range("D3").Value = "'0.05"
$q3Sheet.Range("E3").Value = "'86.27"
$q3Sheet.Range("F3").Value = "'3.97"
$q3Sheet.Range("G3").Value = "'0.0020"
$q3Sheet.Range("H3").Value = 4

# ---------------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: the Q2 row becomes Q3 (value -> 0),
#    the old Q1 row becomes Q2 (value unchanged, 0.01), and a fresh Q1 row
#    is appended below with the data that used to live in row 3.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A3").Copy($total.Range("A4"))

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

# ---------------------------------------------------------------------------
# 4) Restore the original tab selection (the untouched "2022-Q1" sheet was
#    the active tab before this edit) since copying/adding sheets along the
#    way shifts Excel's notion of the "active" tab as a side effect.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
